$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title cell: new wording (note trailing space preserved exactly)
$ws.Range("A1").Value = "BDF Edition 9 "

# Header row tweak: trailing space on "Heure"
$ws.Range("C2").Value = "Heure "

# Row 3 used to hold a sample elimination record (rank, player, time, killer).
# That sample data goes away; only the numeric counter stays, reset to 1.
$ws.Range("B3:D3").Clear()
$ws.Range("A3").Value = 1

# Extend the numbering column down through row 10 (2..8) for the upload list.
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
